$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "65.748.32"
Set-TextValue "E2" "  -5.74%  "
Set-TextValue "D3" "3.276.19"
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "556.26"
Set-TextValue "E5" "  -3.83%  "
Set-TextValue "D6" "183.39"
Set-TextValue "E6" "  -4.69%  "
Set-TextValue "E7" "  -0.08%  "
Set-TextValue "E8" "  -4.02%  "
Set-TextValue "D9" "3.266.75"
Set-TextValue "E9" "  -6.36%  "
Set-TextValue "E10" "  -10.06%  "
Set-TextValue "E11" "  -6.26%  "
Set-TextValue "D12" "47.23"
Set-TextValue "E12" "  -8.29%  "
Set-TextValue "E13" "  -7.15%  "
Set-TextValue "D14" "646.79"
Set-TextValue "E14" "  +0.33%  "
Set-TextValue "E15" "  -5.60%  "
Set-TextValue "D16" "3.795.92"
Set-TextValue "E16" "  -6.37%  "
Set-TextValue "D17" "18.02"
Set-TextValue "E17" "  -1.87%  "
Set-TextValue "D18" "65.758.14"
Set-TextValue "E18" "  -5.69%  "
Set-TextValue "E19" "  -3.33%  "
Set-TextValue "D20" "3.275.89"
Set-TextValue "E20" "  -6.35%  "
Set-TextValue "E21" "  -8.87%  "
Set-TextValue "D22" "0.903"
Set-TextValue "E22" "  -5.00%  "
Set-TextValue "E23" "  +1.45%  "
Set-TextValue "D24" "108.44"
Set-TextValue "E24" "  +9.49%  "
Set-TextValue "E25" "  -8.19%  "
Set-TextValue "D27" "2.68"
Set-TextValue "E27" "  -7.12%  "
Set-TextValue "E28" "  -4.88%  "
Set-TextValue "D29" "8.65"
Set-TextValue "E29" "  -8.03%  "
Set-TextValue "D30" "30.24"
Set-TextValue "E30" "  -7.48%  "
Set-TextValue "D31" "3.94"
Set-TextValue "E31" "  -6.97%  "
Set-TextValue "D32" "6.28"
Set-TextValue "E32" "  -6.80%  "
Set-TextValue "D33" "11.05"
Set-TextValue "E33" "  -5.32%  "
Set-TextValue "E34" "  -4.61%  "
Set-TextValue "D35" "3.765.96"
Set-TextValue "E35" "  +1.61%  "
Set-TextValue "D36" "57.48"
Set-TextValue "E37" "  -0.07%  "
Set-TextValue "D38" "518.88"
Set-TextValue "E38" "  -7.96%  "
Set-TextValue "D39" "3.41"
Set-TextValue "E39" "  -5.65%  "
Set-TextValue "D40" "0.0₃0734"
Set-TextValue "E40" "  -7.26%  "
Set-TextValue "D41" "0.131"
Set-TextValue "E41" "  -1.92%  "
Set-TextValue "D42" "2.72"
Set-TextValue "E42" "  -6.24%  "
Set-TextValue "D43" "3.40"
Set-TextValue "E43" "  -16.91%  "
Set-TextValue "D44" "32.87"
Set-TextValue "E44" "  -4.23%  "
Set-TextValue "E45" "  -10.22%  "
Set-TextValue "E46" "  -6.91%  "
Set-TextValue "D47" "3.20"
Set-TextValue "E47" "  -5.10%  "
Set-TextValue "E48" "  -4.34%  "
Set-TextValue "D49" "2.60"
Set-TextValue "E49" "  -8.61%  "
Set-TextValue "E50" "  +0.13%  "
Set-TextValue "D51" "1.26"
Set-TextValue "E51" "  +1.81%  "
